# The document has three inline pictures whose "Name" (wp:docPr / pic:cNvPr
# name attribute) needs to be swapped:
#   - Footer "default"    (Section 1, Footers(1)) -> Pearson logo: image2.png -> image1.png
#   - Footer "first page" (Section 1, Footers(2)) -> Pearson logo: image2.png -> image1.png
#   - Header "first page" (Section 1, Headers(2)) -> BTec logo   : image1.jpg -> image2.jpg
#
# InlineShape has no settable .Name property in the Word object model, so each
# picture is momentarily converted to a floating Shape (which does expose
# .Name), renamed, then converted back to an inline picture in place.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($inlineShape, $newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# Footer (default) - Pearson logo
$footerDefault = $sec.Footers.Item(1)
Rename-InlinePicture $footerDefault.Range.InlineShapes.Item(1) "image1.png"

# Footer (first page) - Pearson logo
$footerFirst = $sec.Footers.Item(2)
Rename-InlinePicture $footerFirst.Range.InlineShapes.Item(1) "image1.png"

# Header (first page) - BTec logo
$headerFirst = $sec.Headers.Item(2)
Rename-InlinePicture $headerFirst.Range.InlineShapes.Item(1) "image2.jpg"
